$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing I2:I5 values from 4 to 5
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

# Add new row 6 with data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = -5
$ws.Range("H6").Value = 21
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

# Update selection to I6 to match the final state
$ws.Range("I6").Select()
